$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.323.52'
$ws.Range("E2").Value = '  -2.72%  '

$ws.Range("D3").Value = '1.831.06'
$ws.Range("E3").Value = '  -2.57%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '259.51'
$ws.Range("E5").Value = '  -7.65%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").Value = '0.5172'
$ws.Range("E7").Value = '  -2.25%  '

$ws.Range("D8").Value = '0.3229'
$ws.Range("E8").Value = '  -8.26%  '

$ws.Range("D9").Value = '0.06727'
$ws.Range("E9").Value = '  -4.34%  '

$ws.Range("D10").Value = '18.69'
$ws.Range("E10").Value = '  -8.07%  '

$ws.Range("D11").Value = '0.7652'
$ws.Range("E11").Value = '  -6.13%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07679'
$ws.Range("E12").Value = '  -1.63%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.871.80'
$ws.Range("E13").Value = '  -0.42%  '

$ws.Range("D14").Value = '88.64'
$ws.Range("E14").Value = '  -2.07%  '

$ws.Range("D15").Value = '5.012'
$ws.Range("E15").Value = '  -3.64%  '

$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").Value = '14.05'
$ws.Range("E17").Value = '  -3.61%  '

$ws.Range("D19").Value = '0.000007861'
$ws.Range("E19").Value = '  -4.18%  '

$ws.Range("D20").Value = '26.338.36'
$ws.Range("E20").Value = '  -2.81%  '

$ws.Range("D21").Value = '2.084.71'
$ws.Range("E21").Value = '  -1.61%  '

$ws.Range("D22").Value = '4.531'
$ws.Range("E22").Value = '  -4.94%  '

$ws.Range("D23").Value = '9.402'
$ws.Range("E23").Value = '  -7.33%  '

$ws.Range("D24").Value = '5.893'
$ws.Range("E24").Value = '  -5.33%  '

$ws.Range("D25").Value = '2.305'
$ws.Range("E25").Value = '  -3.40%  '

$ws.Range("E26").Value = '  -1.10%  '

$ws.Range("D27").Value = '1.651'
$ws.Range("E27").Value = '  -1.46%  '

$ws.Range("E28").Value = '  -3.98%  '

$ws.Range("D29").Value = '110.71'
$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("D30").Value = '4.176'
$ws.Range("E30").Value = '  -4.92%  '

$ws.Range("D31").Value = '4.108'
$ws.Range("E31").Value = '  -6.06%  '

$ws.Range("D32").Value = '0.08712'
$ws.Range("E32").Value = '  -2.19%  '

$ws.Range("D33").Value = '0.04827'
$ws.Range("E33").Value = '  -1.59%  '

$ws.Range("D34").Value = '1.124'
$ws.Range("E34").Value = '  -4.34%  '

$ws.Range("D35").Value = '2.850'
$ws.Range("E35").Value = '  -1.71%  '

$ws.Range("D36").Value = '0.6785'
$ws.Range("E36").Value = '  -8.92%  '

$ws.Range("D37").Value = '3.095'
$ws.Range("E37").Value = '  -6.49%  '

$ws.Range("D38").Value = '0.01778'
$ws.Range("E38").Value = '  -5.54%  '

$ws.Range("D39").Value = '2.196'
$ws.Range("E39").Value = '  -8.91%  '

$ws.Range("D40").Value = '0.4885'
$ws.Range("E40").Value = '  -8.16%  '

$ws.Range("D41").Value = '0.8972'
$ws.Range("E41").Value = '  -8.05%  '

$ws.Range("D42").Value = '110.85'
$ws.Range("E42").Value = '  -5.36%  '

$ws.Range("D43").Value = '6.132'
$ws.Range("E43").Value = '  -2.91%  '

$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").Value = '7.680'
$ws.Range("E45").Value = '  -6.39%  '

$ws.Range("D46").Value = '0.4189'
$ws.Range("E46").Value = '  -8.93%  '

$ws.Range("D47").Value = '0.1254'
$ws.Range("E47").Value = '  -8.39%  '

$ws.Range("D48").Value = '9.077'
$ws.Range("E48").Value = '  -3.73%  '

$ws.Range("E49").Value = '  -1.20%  '

$ws.Range("D50").Value = '35.26'
$ws.Range("E50").Value = '  -3.93%  '

$ws.Range("D51").Value = '1.413'
$ws.Range("E51").Value = '  -7.54%  '
